# Apply project-status updates to "Apache java projects assignment" workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 99-118 (Keye Li's batch): mark status column C as "Done" now that work
# finished. Row 103 also gets its note (column D) replaced since the earlier
# "no git repo" note was corrected to point at the real failure reason.
for ($r = 99; $r -le 118; $r++) {
    $ws.Range("C$r").Value = "Done"
}
$ws.Range("D103").Value = "Failed to extract data from GetInfo.java"

# Rows 149-158 (Apache Mahout .. Apache Nutch) were previously unassigned;
# assign them to Keye Li and mark them in progress.
for ($r = 149; $r -le 158; $r++) {
    $ws.Range("B$r").Value = "Keye Li"
    $ws.Range("C$r").Value = "in progress"
}

# Reflect the scroll/selection position recorded in the saved workbook.
$ws.Application.ActiveWindow.ScrollRow = 142
$ws.Range("J154").Select()
